# 17.1.2.xlsx update: add the 2020 data point (column Q) to the
# "Proportion of domestic budget funded by domestic taxes" table and
# move the on-screen selection to Q8, matching the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column Q: year header (row 4) and value (row 5) -------------------
# Copy formatting from the existing neighbouring year cells (P4 / P5) so the
# new cells render with the same font/border/alignment as the rest of the
# header and data rows, then overwrite with the new 2020 figures.
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2020

$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 70.3

# --- Selection -------------------------------------------------------------
# The sheet's live selection moves to Q8 (just past the new data).
$ws.Range("Q8").Select()
